# Update cryptos list (GitHub Actions style refresh of coinranking.com scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 43 and 44 swapped position (TheSandbox <-> FraxShare) with updated price/volume ---
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.079"
$ws.Range("E43").Value = "  +8.72%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3706"
$ws.Range("E44").Value = "  +0.97%  "

# --- Price / Volume(1h) refresh for all other rows ---
# (leading "'" forces plain-numeric-looking prices to stay text, matching the
#  original inline-string cell type instead of being auto-converted to a number)
$ws.Range("D2").Value = "25.529.70"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "1.665.40"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'234.38"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("D8").Value = "'0.2579"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.06134"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "1.664.59"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").Value = "'0.06943"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "'4.359"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "'74.93"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "'0.5721"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'0.9996"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "25.529.22"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "'0.000006713"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").Value = "'11.36"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "1.877.91"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "'4.421"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "'8.693"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").Value = "'5.210"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'134.64"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'14.87"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'1.363"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'1.704"
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Value = "'3.949"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "'0.07704"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").Value = "'3.601"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("D33").Value = "'0.04304"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "'2.617"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").Value = "'0.9435"
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("D36").Value = "'0.5983"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "'0.9209"
$ws.Range("E37").Value = "  +11.95%  "
$ws.Range("D38").Value = "'2.480"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "'0.9989"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'103.18"
$ws.Range("E40").Value = "  +5.57%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "'1.823"
$ws.Range("E42").Value = "  +5.71%  "
$ws.Range("D45").Value = "'0.1108"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").Value = "'0.05252"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "'6.112"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "'29.71"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "'7.396"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'0.9971"
$ws.Range("E51").Value = "  -0.18%  "

Write-Host "Cryptos list updated"
